$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 29648
$ws.Cells.Item(2, 2).Value = "Ravi Lucca Abreu"
$ws.Cells.Item(2, 3).Value = "Operacoes"
$ws.Cells.Item(2, 4).Value = "Viagem de negocios"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 45095
$ws.Cells.Item(2, 7).Value = 8681.41

# Row 3
$ws.Cells.Item(3, 1).Value = 72664
$ws.Cells.Item(3, 2).Value = "Dra. Lorena da Mota"
$ws.Cells.Item(3, 3).Value = "TI"
$ws.Cells.Item(3, 5).Value = 8
$ws.Cells.Item(3, 6).Value = 45105
$ws.Cells.Item(3, 7).Value = 8302.549999999999

# Row 4
$ws.Cells.Item(4, 1).Value = 35738
$ws.Cells.Item(4, 2).Value = "Sra. Heloísa Abreu"
$ws.Cells.Item(4, 3).Value = "Engenharia"
$ws.Cells.Item(4, 4).Value = "Consulta medica"
$ws.Cells.Item(4, 6).Value = 45087
$ws.Cells.Item(4, 7).Value = 7276.8

# Row 5
$ws.Cells.Item(5, 1).Value = 58130
$ws.Cells.Item(5, 2).Value = "Davi Moraes"
$ws.Cells.Item(5, 3).Value = "Recursos Humanos"
$ws.Cells.Item(5, 4).Value = "Problemas pessoais"
$ws.Cells.Item(5, 5).Value = 6
$ws.Cells.Item(5, 6).Value = 45101
$ws.Cells.Item(5, 7).Value = 3907.81

# Row 6
$ws.Cells.Item(6, 1).Value = 74404
$ws.Cells.Item(6, 2).Value = "Pedro Henrique Câmara"
$ws.Cells.Item(6, 3).Value = "Operacoes"
$ws.Cells.Item(6, 4).Value = "Problemas pessoais"
$ws.Cells.Item(6, 5).Value = 7
$ws.Cells.Item(6, 6).Value = 45091
$ws.Cells.Item(6, 7).Value = 3389.46

# Row 7
$ws.Cells.Item(7, 1).Value = 90697
$ws.Cells.Item(7, 2).Value = "Sophia da Mata"
$ws.Cells.Item(7, 3).Value = "Atendimento ao Cliente"
$ws.Cells.Item(7, 4).Value = "Problemas pessoais"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 45087
$ws.Cells.Item(7, 7).Value = 9849.07

# Row 8
$ws.Cells.Item(8, 1).Value = 79550
$ws.Cells.Item(8, 2).Value = "Sr. Diego Monteiro"
$ws.Cells.Item(8, 3).Value = "Operacoes"
$ws.Cells.Item(8, 5).Value = 5
$ws.Cells.Item(8, 6).Value = 45103
$ws.Cells.Item(8, 7).Value = 5971.48

# Row 9
$ws.Cells.Item(9, 1).Value = 46353
$ws.Cells.Item(9, 2).Value = "Diego Albuquerque"
$ws.Cells.Item(9, 3).Value = "Financeiro"
$ws.Cells.Item(9, 4).Value = "Viagem de negocios"
$ws.Cells.Item(9, 5).Value = 4
$ws.Cells.Item(9, 6).Value = 45092
$ws.Cells.Item(9, 7).Value = 7676.01

# Row 10
$ws.Cells.Item(10, 1).Value = 89152
$ws.Cells.Item(10, 2).Value = "Caleb Abreu"
$ws.Cells.Item(10, 3).Value = "Recursos Humanos"
$ws.Cells.Item(10, 4).Value = "Doenca"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 45105
$ws.Cells.Item(10, 7).Value = 9657.049999999999

# Row 11
$ws.Cells.Item(11, 1).Value = 14175
$ws.Cells.Item(11, 2).Value = "Maria Eduarda Sousa"
$ws.Cells.Item(11, 3).Value = "Juridico"
$ws.Cells.Item(11, 4).Value = "Problemas pessoais"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 45094
$ws.Cells.Item(11, 7).Value = 5053.13
